$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.905.55"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "3.785.08"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("D7").Value = "3.782.07"
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.447"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.80%  "
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").Value = "4.419.71"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "3.783.74"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "67.877.82"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "457.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.691"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").Value = "3.934.25"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.93%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("E38").Value = "  +6.93%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.16%  "
